# Edit script: insert two new rows (49 and 50) with new price observations
# for "Feria Lagunitas de Puerto Montt - Kiwi", pushing the existing data
# (previously rows 49-126) down to rows 51-128.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows before row 49 -- this shifts rows 49:126 down to 51:128
$ws.Rows("49:50").Insert()

# ---- New row 49 ----
$ws.Cells.Item(49, 1).Value = 4
$ws.Cells.Item(49, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(49, 3).Value = "Los Lagos"
$ws.Cells.Item(49, 4).Value = 44477
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100101
$ws.Cells.Item(49, 8).Value = "Berries"
$ws.Cells.Item(49, 9).Value = 100101007
$ws.Cells.Item(49, 10).Value = "Kiwi"
$ws.Cells.Item(49, 11).Value = "Hayward"
$ws.Cells.Item(49, 12).Value = "Especial"
$ws.Cells.Item(49, 13).Value = 200
$ws.Cells.Item(49, 14).Value = 21000
$ws.Cells.Item(49, 15).Value = 21000
$ws.Cells.Item(49, 16).Value = 21000
$ws.Cells.Item(49, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(49, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(49, 19).Value = 1400
$ws.Cells.Item(49, 20).Value = 15

# ---- New row 50 ----
$ws.Cells.Item(50, 1).Value = 4
$ws.Cells.Item(50, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(50, 3).Value = "Los Lagos"
$ws.Cells.Item(50, 4).Value = 44477
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 10
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100101
$ws.Cells.Item(50, 8).Value = "Berries"
$ws.Cells.Item(50, 9).Value = 100101007
$ws.Cells.Item(50, 10).Value = "Kiwi"
$ws.Cells.Item(50, 11).Value = "Hayward"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 600
$ws.Cells.Item(50, 14).Value = 15000
$ws.Cells.Item(50, 15).Value = 16000
$ws.Cells.Item(50, 16).Value = 15500
$ws.Cells.Item(50, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(50, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(50, 19).Value = 1033
$ws.Cells.Item(50, 20).Value = 15
